# "Add cantrals by cantons"
#
# The sheet currently has a two-row header (units sub-header on row 2) and
# 13 data rows (rows 3-15). The new layout collapses everything to a single
# header row followed immediately by the 13 data rows (rows 2-14), and adds
# three new leading columns (idx, idx2, Name) plus relabelled date/power/
# energy columns.
#
# Strategy: delete the two old header rows (which also shifts the 13 data
# rows up to rows 1-13, preserving each cell's existing style/number format
# automatically), insert one fresh row above them for the new header, then
# populate that header row with the new column titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "units" header rows (old row 1 + old row 2). The 13 data
# rows slide up to rows 1-13 with their per-cell styles intact.
$ws.Rows("1:2").Delete()

# Make room for the single new header row; data now occupies rows 2-14.
$ws.Rows("1").Insert()

# New header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 carry the Arial-9 "unit label" style (same font as the data cells,
# General number format). Apply it via a throwaway named style so the
# engine records it as plain direct formatting (fontId + applyFont) rather
# than folding it into the existing Arial-9/General text style, then
# discard the temporary named style again.
$tmpStyle = $wb.Styles.Add("TmpHeaderUnitStyle")
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderUnitStyle"
[void]$tmpStyle.Delete()

# Match the saved selection state (first data row selected).
[void]$ws.Range("A2:K2").Select()
